$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits (per sharedStrings diff) ---
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Row height changes (18.75 -> 19.5) ---
$ws.Range("A1:L3").EntireRow.RowHeight = 19.5

# --- Font color normalization for the Pincode/Phone columns (theme -> explicit black) ---
$ws.Range("I2:I3").Font.Color = 0
$ws.Range("K2:K3").Font.Color = 0
